$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 345.8
$ws.Range("I12").Value = 320
$ws.Range("K12").Value = 320
$ws.Range("M12").Value = -150
$ws.Range("H17").Value = 270583.06
$ws.Range("J17").Value = 278029.8
$ws.Range("L17").Value = 834089.3999999999
$ws.Range("N17").Value = -834425.3999999999
$ws.Range("H18").Value = 486
$ws.Range("H40").Value = 3795.111
$ws.Range("I40").Value = 8329.666999999999
$ws.Range("J40").Value = 1527.8334
$ws.Range("K40").Value = 8329.666999999999
$ws.Range("L40").Value = 1527.8334
$ws.Range("M40").Value = -8154.666999999999
$ws.Range("N40").Value = -1877.8334
$ws.Range("H48").Value = 1673
$ws.Range("I48").Value = 1450
$ws.Range("J48").Value = 1784.5
$ws.Range("K48").Value = 4350
$ws.Range("L48").Value = 5353.5
$ws.Range("M48").Value = -4058
$ws.Range("N48").Value = -5937.5
$ws.Range("H51").Value = 14740.689
$ws.Range("I51").Value = 21397.4
$ws.Range("J51").Value = 13353.875
$ws.Range("K51").Value = 21397.4
$ws.Range("L51").Value = 13353.875
$ws.Range("M51").Value = -20913.4
$ws.Range("N51").Value = -14321.875
$ws.Range("H56").Value = 1673
$ws.Range("I56").Value = 1450
$ws.Range("J56").Value = 1784.5
$ws.Range("K56").Value = 4350
$ws.Range("L56").Value = 5353.5
$ws.Range("M56").Value = -3816
$ws.Range("N56").Value = -6421.5
$ws.Range("H58").Value = 433.93332
$ws.Range("I58").Value = 433.93332
$ws.Range("K58").Value = 1301.79996
$ws.Range("M58").Value = -1151.79996
$ws.Range("H74").Value = 4326.857
$ws.Range("I74").Value = 4804.75
$ws.Range("K74").Value = 4804.75
$ws.Range("M74").Value = -3868.75
$ws.Range("H76").Value = 3973.4614
$ws.Range("I76").Value = 3304.3333
$ws.Range("J76").Value = 4547
$ws.Range("K76").Value = 3304.3333
$ws.Range("L76").Value = 4547
$ws.Range("M76").Value = -2989.3333
$ws.Range("N76").Value = -5177
$ws.Range("H77").Value = 4326.857
$ws.Range("I77").Value = 4804.75
$ws.Range("K77").Value = 24023.75
$ws.Range("M77").Value = -19343.75
$ws.Range("H79").Value = 3973.4614
$ws.Range("I79").Value = 3304.3333
$ws.Range("J79").Value = 4547
$ws.Range("K79").Value = 3304.3333
$ws.Range("L79").Value = 4547
$ws.Range("M79").Value = -2212.3333
$ws.Range("N79").Value = -6731
$ws.Range("H88").Value = 3966.6667
$ws.Range("I88").Value = 3450
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 3450
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -3044
$ws.Range("N88").Value = -5812
$ws.Range("H91").Value = 3966.6667
$ws.Range("I91").Value = 3450
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 3450
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = -2046
$ws.Range("N91").Value = -7808
$ws.Range("H135").Value = 1025
$ws.Range("J135").Value = 725
$ws.Range("L135").Value = 6525
$ws.Range("N135").Value = -11595
$ws.Range("H137").Value = 61197.65
$ws.Range("I137").Value = 125354.445
$ws.Range("K137").Value = 376063.335
$ws.Range("M137").Value = -373513.335
$ws.Range("H138").Value = 2547.2104
$ws.Range("I138").Value = 1304.4814
$ws.Range("K138").Value = 3913.4442
$ws.Range("M138").Value = 1226.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6542.493
$ws.Range("I32").Value = 4058.8462
$ws.Range("J32").Value = 33448.668
$ws.Range("K32").Value = 4058.8462
$ws.Range("L32").Value = 33448.668
$ws.Range("M32").Value = -3771.8462
$ws.Range("N32").Value = -34022.668
$ws.Range("H61").Value = 4170.45
$ws.Range("I61").Value = 2529.9412
$ws.Range("K61").Value = 2529.9412
$ws.Range("M61").Value = -2317.9412
$ws.Range("H132").Value = 1898.0256
$ws.Range("I132").Value = 1769.0667
$ws.Range("K132").Value = 5307.2001
$ws.Range("M132").Value = -2777.2001
$ws.Range("H136").Value = 4170.45
$ws.Range("I136").Value = 2529.9412
$ws.Range("K136").Value = 7589.823600000001
$ws.Range("M136").Value = -5039.823600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2639.5
$ws.Range("I105").Value = 2629.4167
$ws.Range("J105").Value = 2700
$ws.Range("K105").Value = 2629.4167
$ws.Range("L105").Value = 2700
$ws.Range("M105").Value = -882.4167000000002
$ws.Range("N105").Value = -6194

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 9799
$ws.Range("I103").Value = 9799
$ws.Range("K103").Value = 9799
$ws.Range("M103").Value = -8627
$ws.Range("H105").Value = 4166.9585
$ws.Range("I105").Value = 1373.3914
$ws.Range("K105").Value = 1373.3914
$ws.Range("M105").Value = 373.6086
$ws.Range("H132").Value = 3048.913
$ws.Range("I132").Value = 2274.9333
$ws.Range("K132").Value = 6824.7999
$ws.Range("M132").Value = -4294.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 676.53845
$ws.Range("I113").Value = 396.4
$ws.Range("K113").Value = 1189.2
$ws.Range("M113").Value = 980.8000000000002
$ws.Range("H131").Value = 10001526
$ws.Range("I131").Value = 250000750
$ws.Range("J131").Value = 1558.5209
$ws.Range("K131").Value = 750002250
$ws.Range("L131").Value = 4675.5627
$ws.Range("M131").Value = -749997210
$ws.Range("N131").Value = -14755.5627

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 49998
$ws.Range("J63").Value = 49998
$ws.Range("L63").Value = 49998
$ws.Range("N63").Value = -51370
$ws.Range("H66").Value = 49998
$ws.Range("J66").Value = 49998
$ws.Range("L66").Value = 149994
$ws.Range("N66").Value = -156858
$ws.Range("H104").Value = 73175.336
$ws.Range("J104").Value = 73175.336
$ws.Range("L104").Value = 73175.336
$ws.Range("N104").Value = -80163.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1093.1111
$ws.Range("I22").Value = 837
$ws.Range("J22").Value = 1237.8695
$ws.Range("K22").Value = 837
$ws.Range("L22").Value = 1237.8695
$ws.Range("M22").Value = -542
$ws.Range("N22").Value = -1827.8695
$ws.Range("H27").Value = 1093.1111
$ws.Range("I27").Value = 837
$ws.Range("J27").Value = 1237.8695
$ws.Range("K27").Value = 837
$ws.Range("L27").Value = 1237.8695
$ws.Range("M27").Value = -730
$ws.Range("N27").Value = -1451.8695
$ws.Range("H46").Value = 3206.5557
$ws.Range("I46").Value = 2551.1428
$ws.Range("J46").Value = 5500.5
$ws.Range("K46").Value = 2551.1428
$ws.Range("L46").Value = 5500.5
$ws.Range("M46").Value = -2363.1428
$ws.Range("N46").Value = -5876.5
$ws.Range("H57").Value = 15045
$ws.Range("J57").Value = 15045
$ws.Range("L57").Value = 15045
$ws.Range("N57").Value = -16177
$ws.Range("H132").Value = 8623.875
$ws.Range("I132").Value = 9118.799999999999
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 27356.4
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -24826.4
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 18364.666
$ws.Range("I39").Value = 18364.666
$ws.Range("K39").Value = 18364.666
$ws.Range("M39").Value = -17951.666
$ws.Range("H42").Value = 49997.5
$ws.Range("I42").Value = 49997.5
$ws.Range("K42").Value = 49997.5
$ws.Range("M42").Value = -49619.5
$ws.Range("H132").Value = 1059.3226
$ws.Range("I132").Value = 1059.4814
$ws.Range("J132").Value = 1058.25
$ws.Range("K132").Value = 3178.4442
$ws.Range("L132").Value = 3174.75
$ws.Range("M132").Value = -648.4441999999999
$ws.Range("N132").Value = -8234.75
